# Apply refreshed monthly ADR data to Sheet1 (rows shift from Jun-Jan to Aug-Jul window)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object "object[,]" 11,32

# row 2: Aug
$data[0,0] = "Aug"
$data[0,1] = 0
$data[0,2] = 0
$data[0,3] = 0
$data[0,4] = 0
$data[0,5] = 0
$data[0,6] = 0
$data[0,7] = 0
$data[0,8] = 41.03178571428571
$data[0,9] = 43.20724137931035
$data[0,10] = 43.35515151515152
$data[0,11] = 67.57973684210526
$data[0,12] = 64.7193023255814
$data[0,13] = 45.30756756756757
$data[0,14] = 41.7053125
$data[0,15] = 41.0888
$data[0,16] = 41.54130434782609
$data[0,17] = 46.15476190476191
$data[0,18] = 53.79944444444445
$data[0,19] = 52.6635
$data[0,20] = 42.33071428571429
$data[0,21] = 37.64
$data[0,22] = 37.71541666666666
$data[0,23] = 37.19391304347826
$data[0,24] = 37.671
$data[0,25] = 54.76111111111111
$data[0,26] = 60.02206896551724
$data[0,27] = 58.6721052631579
$data[0,28] = 61.83333333333334
$data[0,29] = 61.21875
$data[0,30] = 47.82727272727273
$data[0,31] = 52.17368421052631

# row 3: Sep
$data[1,0] = "Sep"
$data[1,1] = 61.82958333333334
$data[1,2] = 60.62612903225807
$data[1,3] = 52.075
$data[1,4] = 44.89909090909091
$data[1,5] = 43.16090909090909
$data[1,6] = 44.332
$data[1,7] = 45.7125
$data[1,8] = 55.50894736842105
$data[1,9] = 49.98708333333334
$data[1,10] = 48.82761904761905
$data[1,11] = 46.9078947368421
$data[1,12] = 45.17173913043479
$data[1,13] = 45.14346153846154
$data[1,14] = 41.001
$data[1,15] = 57.71897435897436
$data[1,16] = 56.03558823529412
$data[1,17] = 40.24842105263158
$data[1,18] = 46.1605
$data[1,19] = 49.85333333333334
$data[1,20] = 49.1325
$data[1,21] = 46.5945
$data[1,22] = 51.443
$data[1,23] = 46.95647058823529
$data[1,24] = 40.33214285714286
$data[1,25] = 35.83461538461539
$data[1,26] = 36.08846153846154
$data[1,27] = 43.75
$data[1,28] = 47.32333333333334
$data[1,29] = 73.20588235294117
$data[1,30] = 66.01590909090909
$data[1,31] = 0

# row 4: Oct
$data[2,0] = "Oct"
$data[2,1] = 116
$data[2,2] = 124.5
$data[2,3] = 124.5
$data[2,4] = 41
$data[2,5] = 73
$data[2,6] = 104.2173913043478
$data[2,7] = 105.2727272727273
$data[2,8] = 70.5
$data[2,9] = 39.66666666666666
$data[2,10] = 39.66666666666666
$data[2,11] = 39.66666666666666
$data[2,12] = 87.46153846153847
$data[2,13] = 109.0869565217391
$data[2,14] = 109.0869565217391
$data[2,15] = 75.52380952380952
$data[2,16] = 53.42857142857143
$data[2,17] = 59.75
$data[2,18] = 53.5
$data[2,19] = 75.30434782608695
$data[2,20] = 96.42307692307692
$data[2,21] = 104.84
$data[2,22] = 76.47619047619048
$data[2,23] = 63.5
$data[2,24] = 41
$data[2,25] = 49.25
$data[2,26] = 65.93333333333334
$data[2,27] = 87
$data[2,28] = 87
$data[2,29] = 58.66666666666666
$data[2,30] = 41
$data[2,31] = 41

# row 5: Nov
$data[3,0] = "Nov"
$data[3,1] = 41
$data[3,2] = 62
$data[3,3] = 81
$data[3,4] = 81.2
$data[3,5] = 52
$data[3,6] = 41
$data[3,7] = 41
$data[3,8] = 38
$data[3,9] = 49.5
$data[3,10] = 41
$data[3,11] = 41
$data[3,12] = 37
$data[3,13] = 35
$data[3,14] = 35
$data[3,15] = 0
$data[3,16] = 37.775
$data[3,17] = 45.6
$data[3,18] = 47.16666666666666
$data[3,19] = 34.53333333333333
$data[3,20] = 29.84
$data[3,21] = 33.33333333333334
$data[3,22] = 35.5
$data[3,23] = 37.3925
$data[3,24] = 42.89000000000001
$data[3,25] = 43.2175
$data[3,26] = 37.95666666666667
$data[3,27] = 38.08
$data[3,28] = 37
$data[3,29] = 35.5
$data[3,30] = 37.35
$data[3,31] = 0

# row 6: Dec
$data[4,0] = "Dec"
$data[4,1] = 45.47333333333334
$data[4,2] = 44
$data[4,3] = 39
$data[4,4] = 37
$data[4,5] = 68
$data[4,6] = 68
$data[4,7] = 71.40000000000001
$data[4,8] = 66
$data[4,9] = 66
$data[4,10] = 39
$data[4,11] = 26.78333333333333
$data[4,12] = 26.78333333333333
$data[4,13] = 24.74
$data[4,14] = 39
$data[4,15] = 0
$data[4,16] = 0
$data[4,17] = 0
$data[4,18] = 64.59999999999999
$data[4,19] = 64.59999999999999
$data[4,20] = 39.26
$data[4,21] = 39.125
$data[4,22] = 43.125
$data[4,23] = 47
$data[4,24] = 39
$data[4,25] = 37
$data[4,26] = 37
$data[4,27] = 37
$data[4,28] = 0
$data[4,29] = 47
$data[4,30] = 47
$data[4,31] = 37.8

# row 7: Jan
$data[5,0] = "Jan"
$data[5,1] = 36.8
$data[5,2] = 0
$data[5,3] = 0
$data[5,4] = 0
$data[5,5] = 0
$data[5,6] = 0
$data[5,7] = 0
$data[5,8] = 0
$data[5,9] = 0
$data[5,10] = 0
$data[5,11] = 0
$data[5,12] = 0
$data[5,13] = 0
$data[5,14] = 0
$data[5,15] = 0
$data[5,16] = 0
$data[5,17] = 0
$data[5,18] = 0
$data[5,19] = 66.33333333333333
$data[5,20] = 66.33333333333333
$data[5,21] = 53
$data[5,22] = 0
$data[5,23] = 37
$data[5,24] = 0
$data[5,25] = 0
$data[5,26] = 0
$data[5,27] = 0
$data[5,28] = 0
$data[5,29] = 0
$data[5,30] = 0
$data[5,31] = 0

# row 8: Feb
$data[6,0] = "Feb"
$data[6,1] = 0
$data[6,2] = 0
$data[6,3] = 0
$data[6,4] = 0
$data[6,5] = 0
$data[6,6] = 0
$data[6,7] = 0
$data[6,8] = 34
$data[6,9] = 37.4
$data[6,10] = 37.4
$data[6,11] = 0
$data[6,12] = 0
$data[6,13] = 0
$data[6,14] = 0
$data[6,15] = 74.5
$data[6,16] = 0
$data[6,17] = 0
$data[6,18] = 0
$data[6,19] = 0
$data[6,20] = 0
$data[6,21] = 0
$data[6,22] = 0
$data[6,23] = 0
$data[6,24] = 0
$data[6,25] = 0
$data[6,26] = 0
$data[6,27] = 0
$data[6,28] = 0
$data[6,29] = 0
$data[6,30] = 0
$data[6,31] = 0

# row 9: Mar
$data[7,0] = "Mar"
$data[7,1] = 0
$data[7,2] = 0
$data[7,3] = 149
$data[7,4] = 129
$data[7,5] = 189
$data[7,6] = 121.5
$data[7,7] = 69
$data[7,8] = 87.33333333333333
$data[7,9] = 87.33333333333333
$data[7,10] = 62.33333333333334
$data[7,11] = 47.75
$data[7,12] = 110
$data[7,13] = 110
$data[7,14] = 143
$data[7,15] = 186
$data[7,16] = 81.5
$data[7,17] = 49
$data[7,18] = 0
$data[7,19] = 0
$data[7,20] = 0
$data[7,21] = 0
$data[7,22] = 0
$data[7,23] = 0
$data[7,24] = 0
$data[7,25] = 0
$data[7,26] = 0
$data[7,27] = 0
$data[7,28] = 0
$data[7,29] = 0
$data[7,30] = 52
$data[7,31] = 41

# row 10: Apr
$data[8,0] = "Apr"
$data[8,1] = 39.6
$data[8,2] = 39.6
$data[8,3] = 39.5
$data[8,4] = 52.225
$data[8,5] = 76.6970588235294
$data[8,6] = 82.5242105263158
$data[8,7] = 72.56
$data[8,8] = 61.6375
$data[8,9] = 61.47692307692308
$data[8,10] = 102.6666666666667
$data[8,11] = 0
$data[8,12] = 0
$data[8,13] = 0
$data[8,14] = 0
$data[8,15] = 0
$data[8,16] = 0
$data[8,17] = 0
$data[8,18] = 0
$data[8,19] = 0
$data[8,20] = 0
$data[8,21] = 0
$data[8,22] = 0
$data[8,23] = 0
$data[8,24] = 0
$data[8,25] = 0
$data[8,26] = 0
$data[8,27] = 0
$data[8,28] = 0
$data[8,29] = 0
$data[8,30] = 0
$data[8,31] = 0

# row 11: Jun
$data[9,0] = "Jun"
$data[9,1] = 0
$data[9,2] = 0
$data[9,3] = 0
$data[9,4] = 0
$data[9,5] = 0
$data[9,6] = 0
$data[9,7] = 0
$data[9,8] = 0
$data[9,9] = 0
$data[9,10] = 0
$data[9,11] = 0
$data[9,12] = 43
$data[9,13] = 0
$data[9,14] = 0
$data[9,15] = 0
$data[9,16] = 0
$data[9,17] = 0
$data[9,18] = 0
$data[9,19] = 0
$data[9,20] = 0
$data[9,21] = 0
$data[9,22] = 0
$data[9,23] = 0
$data[9,24] = 0
$data[9,25] = 0
$data[9,26] = 0
$data[9,27] = 0
$data[9,28] = 0
$data[9,29] = 0
$data[9,30] = 0
$data[9,31] = 0

# row 12: Jul
$data[10,0] = "Jul"
$data[10,1] = 0
$data[10,2] = 0
$data[10,3] = 0
$data[10,4] = 0
$data[10,5] = 0
$data[10,6] = 0
$data[10,7] = 0
$data[10,8] = 0
$data[10,9] = 0
$data[10,10] = 0
$data[10,11] = 0
$data[10,12] = 0
$data[10,13] = 0
$data[10,14] = 0
$data[10,15] = 0
$data[10,16] = 0
$data[10,17] = 0
$data[10,18] = 0
$data[10,19] = 0
$data[10,20] = 0
$data[10,21] = 0
$data[10,22] = 0
$data[10,23] = 0
$data[10,24] = 0
$data[10,25] = 89.5
$data[10,26] = 109.5
$data[10,27] = 109.5
$data[10,28] = 89.5
$data[10,29] = 0
$data[10,30] = 0
$data[10,31] = 0

$ws.Range("A2:AF12").Value = $data

# Carry over the header-row formatting (style, bold, border, centered) to the newly added rows 10-12
$ws.Range("A9").Copy() | Out-Null
$ws.Range("A10:A12").PasteSpecial(-4122)
$excel.CutCopyMode = 0
